# Upgrade config file handling
# Append a new trailing data row (row 76) to each of the four sensor-log
# worksheets, mirroring the existing row layout (time / length / ID /
# actual length / checksum + their *_DEC numeric companions).

$wb = $excel.ActiveWorkbook

$newRows = @{
    "ROW35-FE-LIFTER" = @{
        A = "2025-03-07 11:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    }
    "ROW35-MID-LIFTER" = @{
        A = "2025-03-07 11:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    }
    "ROW02-FE-LIFTER" = @{
        A = "2025-03-07 11:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    }
    "ROW02-MID-LIFTER" = @{
        A = "2025-03-07 11:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
}

foreach ($sheetName in $newRows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = $newRows[$sheetName]

    $targetRow = 76

    $ws.Range("A$targetRow").Value = $row.A
    $ws.Range("B$targetRow").Value = $row.B
    $ws.Range("C$targetRow").Value = $row.C
    $ws.Range("D$targetRow").Value = $row.D
    $ws.Range("E$targetRow").Value = $row.E

    $ws.Range("F$targetRow").Value = $row.F

    # G holds a 24-digit numeric-looking identifier. Plain assignment would
    # coerce it through a double and lose precision, so force the cell to
    # Text first (the normal Excel way to keep a long digit string intact).
    $ws.Range("G$targetRow").NumberFormat = "@"
    $ws.Range("G$targetRow").Value = $row.G

    $ws.Range("H$targetRow").Value = $row.H
    $ws.Range("I$targetRow").Value = $row.I
}
